$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect it so the cells below can be updated.
$ws.Unprotect()

# Update the confidentiality note date from 2021-04-08 to 2021-04-09
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-09 for illustrative purposes only and are subject to change."

# Update D2:E13 values
$ws.Range("D2").Value = 0.03012389620326952
$ws.Range("E2").Value = -0.001780626780626737

$ws.Range("D3").Value = 0.02378853871524316
$ws.Range("E3").Value = 0.000476303881876472

$ws.Range("D4").Value = 0.05176143405521942
$ws.Range("E4").Value = -0.002850356294536782

$ws.Range("D5").Value = 0.1342296861503509
$ws.Range("E5").Value = 0.008825526137135054

$ws.Range("D6").Value = 0.02957327855443099
$ws.Range("E6").Value = -0.00838414634146345

$ws.Range("D7").Value = 0.1210429479355334
$ws.Range("E7").Value = 0.0101176173011257

$ws.Range("D8").Value = 0.1004052729113233
$ws.Range("E8").Value = 0.009484066767829979

$ws.Range("D9").Value = 0.02734380744389739
$ws.Range("E9").Value = 0.007643884892086339

$ws.Range("D10").Value = 0.1205380141940744
$ws.Range("E10").Value = 0.007200000000000095

$ws.Range("D11").Value = 0.2562758719780215
$ws.Range("E11").Value = 0.008936550491510209

$ws.Range("D12").Value = 0.1049172518586359
$ws.Range("E12").Value = 0.0009650646593322154

$ws.Range("D13").Value = 0.9999999999999998
$ws.Range("E13").Value = 0.006392130892711734

# Restore sheet protection (the sheet was protected before editing)
$ws.Protect()
